$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.883.53'
$ws.Range('E2').Value = '  +1.30%  '

$ws.Range('D3').Value = '1.844.45'
$ws.Range('E3').Value = '  +1.67%  '

$ws.Range('E4').Value = '  +0.28%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.16'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.01%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.005'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.29%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4706'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3670'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07155'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.79%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9295'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.11%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.58'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.32%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07696'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.79%  '

$ws.Range('D13').Value = '1.829.86'
$ws.Range('E13').Value = '  +0.82%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.281'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.32%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.396'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.50%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.24'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.20%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.40%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008626'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.26%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.26%  '

$ws.Range('D20').Value = '26.934.54'
$ws.Range('E20').Value = '  +1.36%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.46'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.22%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.018'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.28%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.61'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.06%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.932'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.93%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.87'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.04%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.24'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.54%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.012'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.15%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.40'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.89%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.879'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.64%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08850'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.54%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.212'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.75%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.180'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.08%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7472'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.34%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.778'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.48%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.473'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.15%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.083'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.85%  '

$ws.Range('E37').Value = '  +0.20%  '

$ws.Range('E38').Value = '  +2.19%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.959'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.38%  '

$ws.Range('E40').Value = '  +2.04%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.981'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.16%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1511'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.07%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.158'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.52%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.52'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.47%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4708'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.13%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.006'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.37%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.98'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.65%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.597'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.79%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '65.60'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.75%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06043'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.09%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8932'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.23%  '
